# Applies the cryptos.xlsx price/volume/ranking refresh described in the commit diff.
# (GitHub Actions scheduled data refresh: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.224.80"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "2.095.88"
$ws.Range("E3").Value = "  +2.82%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'229.72"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").Value = "'0.614"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").Value = "'60.83"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("D10").Value = "'0.0843"
$ws.Range("E10").Value = "  +2.79%  "
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").Value = "2.406.66"
$ws.Range("E12").Value = "  +2.76%  "
$ws.Range("D13").Value = "'22.45"
$ws.Range("E13").Value = "  +5.18%  "
$ws.Range("D14").Value = "'14.66"
$ws.Range("E15").Value = "  +6.41%  "
$ws.Range("D16").Value = "'0.773"
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("D17").Value = "2.097.71"
$ws.Range("E17").Value = "  +2.80%  "
$ws.Range("D18").Value = "38.156.21"
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("B19").Value = "Litecoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D19").Value = "'70.36"
$ws.Range("E19").Value = "  +0.69%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'6.00"
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("E21").Value = "  +1.26%  "
$ws.Range("D22").Value = "'224.51"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("E25").Value = "  +3.39%  "
$ws.Range("D26").Value = "'170.25"
$ws.Range("E26").Value = "  +1.79%  "
$ws.Range("D27").Value = "'9.43"
$ws.Range("E27").Value = "  +0.86%  "
$ws.Range("D28").Value = "'0.131"
$ws.Range("E28").Value = "  +1.16%  "
$ws.Range("D29").Value = "'19.03"
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("E30").Value = "  +5.27%  "
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").Value = "'2.38"
$ws.Range("E32").Value = "  +8.99%  "
$ws.Range("D33").Value = "'4.70"
$ws.Range("E33").Value = "  +3.60%  "
$ws.Range("D34").Value = "'4.43"
$ws.Range("E34").Value = "  +0.37%  "
$ws.Range("D35").Value = "'0.0606"
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("D36").Value = "'6.53"
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("D37").Value = "'2.40"
$ws.Range("E37").Value = "  +4.98%  "
$ws.Range("D38").Value = "'3.52"
$ws.Range("E38").Value = "  +6.98%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("E40").Value = "  +1.51%  "
$ws.Range("D41").Value = "1.548.96"
$ws.Range("E41").Value = "  +1.42%  "
$ws.Range("D42").Value = "'100.07"
$ws.Range("E42").Value = "  +3.85%  "
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "'0.0905"
$ws.Range("E45").Value = "  -1.19%  "
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("E47").Value = "  +1.78%  "
$ws.Range("E48").Value = "  +1.65%  "
$ws.Range("D49").Value = "'7.25"
$ws.Range("E49").Value = "  +1.90%  "
$ws.Range("D50").Value = "'3.00"
$ws.Range("E50").Value = "  +0.83%  "
$ws.Range("D51").Value = "2.294.62"
$ws.Range("E51").Value = "  +2.86%  "
